$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I0 and IF headers with the same style as the other headers ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# --- Data rows 2..31: I = 1 (constant), J = same value as H ---
$lastRow = 31
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}

$ws.Range("A1").Select() | Out-Null
